$d = $word.ActiveDocument

# --- Edit 1: "Prompt: ... varying heights." -> "... varying heights using Perlin noise." ---
$rng1 = $d.Content
$rng1.Find.Execute("varying heights.", $false, $false, $false, $false, $false, $true, 1, $false, "varying heights using Perlin noise.", 2)

# --- Edit 2: "PerlinNoise = Random.Range(0,5) " -> "PerlinNoise *= Random.Range(0,5) " ---
$rng2 = $d.Content
$rng2.Find.Execute("PerlinNoise = Random.Range(0,5)", $false, $false, $false, $false, $false, $true, 1, $false, "PerlinNoise *= Random.Range(0,5)", 2)

# --- Edit 3: the "Real Information" row's actual start time "5:00PM" -> "5:22PM" ---
$rng3 = $d.Content
$rng3.Find.Execute("Real Information", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Collapse(0)
$rng3.MoveEnd(1, 200)
$sub = $rng3.Duplicate
$sub.Find.Execute("5:00", $false, $false, $false, $false, $false, $true, 1, $false, "5:22", 2)
